# Automatic update of files.
#
# Applies the diff between the previous and current snapshot of the
# "HEBY" logging overview sheet:
#   1. The "Förändrad" (changed) date in column C bumps by one day
#      (45183 -> 45184) for every data row (rows 2-28).
#   2. The link columns (S:Y) for the rows touched by this update get a
#      second HYPERLINK() argument with the case id as display text.
#      Column S ends up with a malformed argument (matches the faulty
#      text produced upstream), while T:Y receive a well-formed second
#      argument. Column Y additionally switches from a literal text
#      value (which used ";" as argument separator and was missing the
#      leading "=", so it never was an actual formula) to a real
#      HYPERLINK formula using ",".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the "changed" date for every row (2-28) from 45183 to 45184.
$ws.Range("C2:C28").Value = 45184

# 2. Update the link formulas for the affected rows (2, 3, 4, 5, 15).

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 31572-2023.xlsx, "A 31572-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 31572-2023.png", "A 31572-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/knärot/A 31572-2023.png", "A 31572-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 32292-2023.xlsx, "A 32292-2023"")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 32292-2023.png", "A 32292-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 32299-2023.xlsx, "A 32299-2023"")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 32299-2023.png", "A 32299-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 32785-2023.xlsx, "A 32785-2023"")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 32785-2023.png", "A 32785-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("U15").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/knärot/A 31569-2023.png", "A 31569-2023")'
$ws.Range("V15").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 31569-2023.docx", "A 31569-2023")'
$ws.Range("W15").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 31569-2023.docx", "A 31569-2023")'
$ws.Range("X15").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 31569-2023.docx", "A 31569-2023")'
$ws.Range("Y15").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 31569-2023.docx", "A 31569-2023")'
